$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projektplan")

# --- Column width / visibility adjustments ---
$ws.Columns("B").ColumnWidth = 41.333333333333336
$ws.Columns("C").ColumnWidth = 22.666666666666668
$ws.Columns("D").ColumnWidth = 1.3333333333333333
$ws.Range("I1:M1").EntireColumn.Hidden = $true
$ws.Columns("T").ColumnWidth = 3.1666666666666665

# --- Project start date formula: one day earlier (Heute-3 -> Heute-4) ---
$ws.Range("E3").FormulaArray = "=Heute-4"

# --- Row 9 ("Schnittstellen zur RKI API aufbauen"): clear assignee + start/end dates ---
$ws.Range("C9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Rows("9").RowHeight = 20.1

# --- Row 11 ("Prognose"): clear assignee + start/end dates ---
$ws.Range("C11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("F11").ClearContents()

# --- Row 13 ("Entwicklung Prognosefunktion"): shift start/end one day earlier ---
$ws.Range("E13").Value2 = 44242
$ws.Range("F13").Value2 = 44243

# --- Row 15 ("Rest"): clear assignee + start/end dates ---
$ws.Range("C15").ClearContents()
$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Rows("15").RowHeight = 20.1

# --- Row 16 ("Entwurf Rest-Design"): update assignee list ---
$ws.Range("C16").Value = "Benjamin, Nico, Thomas"

# --- Row 17 ("Entwicklung Rest-Design"): end date moves one day earlier ---
$ws.Range("F17").Value2 = 44243

# --- Row 18 ("Schnittstellen zur RKI API aufbauen"): end date moves one day earlier ---
$ws.Range("F18").Value2 = 44243

# --- Row 19 ("Entwicklung Vorlagen JSON Representation"): start/end updated ---
$ws.Range("E19").Value2 = 44243
$ws.Range("F19").Value2 = 44243

# --- Row 20 ("Dokumentation / Organisatorisches" footer): clear assignee + start/end dates ---
$ws.Range("C20").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Rows("20").RowHeight = 20.1

# --- Selection cursor moved ---
$ws.Range("W19").Select()
